$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (Price-or-$null, Volume).
# Price column is D, Volume(1h) column is E.
# Price = $null means that row's Price cell is unchanged by this update.
#
# Some of the new Price strings look like plain decimal numbers (e.g.
# "567.00", "141.51"); Excel's normal cell-value coercion would silently
# turn those into numeric values (losing the trailing zero / exact text).
# To keep them as literal text - matching how the sheet already stores
# every Price cell as text - those particular cells are pre-formatted as
# Text ("@") before the value is assigned.

$updates = @(
    @{ Row = 2;  Price = "59.937.90"; Volume = "  +4.07%  " },
    @{ Row = 3;  Price = "3.020.99";  Volume = "  +2.76%  " },
    @{ Row = 4;  Price = $null;       Volume = "  +0.11%  " },
    @{ Row = 5;  Price = "567.00";    Volume = "  +3.27%  " },
    @{ Row = 6;  Price = "141.51";    Volume = "  +8.50%  " },
    @{ Row = 7;  Price = $null;       Volume = "  -0.06%  " },
    @{ Row = 8;  Price = $null;       Volume = "  +2.22%  " },
    @{ Row = 9;  Price = "3.010.53";  Volume = "  +2.58%  " },
    @{ Row = 10; Price = $null;       Volume = "  +6.82%  " },
    @{ Row = 11; Price = "5.38";      Volume = "  +13.23%  " },
    @{ Row = 12; Price = $null;       Volume = "  +3.42%  " },
    @{ Row = 13; Price = $null;       Volume = "  +5.47%  " },
    @{ Row = 14; Price = "34.16";     Volume = "  +3.85%  " },
    @{ Row = 15; Price = $null;       Volume = "  +0.71%  " },
    @{ Row = 16; Price = "3.519.85";  Volume = "  +2.83%  " },
    @{ Row = 17; Price = "7.19";      Volume = "  +4.34%  " },
    @{ Row = 18; Price = "3.018.13";  Volume = "  +2.82%  " },
    @{ Row = 19; Price = "59.874.37"; Volume = "  +3.93%  " },
    @{ Row = 20; Price = "440.36";    Volume = "  +5.57%  " },
    @{ Row = 21; Price = "13.70";     Volume = "  +3.45%  " },
    @{ Row = 22; Price = "0.721";     Volume = "  +4.92%  " },
    @{ Row = 23; Price = $null;       Volume = "  +2.09%  " },
    @{ Row = 24; Price = "13.42";     Volume = "  +2.49%  " },
    @{ Row = 25; Price = "80.82";     Volume = "  +1.28%  " },
    @{ Row = 26; Price = $null;       Volume = "  -0.05%  " },
    @{ Row = 27; Price = "2.23";      Volume = "  +12.47%  " },
    @{ Row = 28; Price = $null;       Volume = "  +0.15%  " },
    @{ Row = 29; Price = $null;       Volume = "  +4.00%  " },
    @{ Row = 30; Price = "7.90";      Volume = "  +6.00%  " },
    @{ Row = 31; Price = $null;       Volume = "  +5.90%  " },
    @{ Row = 32; Price = "26.10";     Volume = "  +3.50%  " },
    @{ Row = 33; Price = "0.107";     Volume = "  +10.90%  " },
    @{ Row = 34; Price = "0.0₃0798";  Volume = "  +16.12%  " },
    @{ Row = 35; Price = $null;       Volume = "  +7.12%  " },
    @{ Row = 36; Price = "5.96";      Volume = "  +5.03%  " },
    @{ Row = 37; Price = "2.13";      Volume = "  +2.53%  " },
    @{ Row = 38; Price = "49.44";     Volume = "  +2.47%  " },
    @{ Row = 39; Price = "8.64";      Volume = "  -0.51%  " },
    @{ Row = 40; Price = $null;       Volume = "  +10.96%  " },
    @{ Row = 41; Price = "406.97";    Volume = "  +8.04%  " },
    @{ Row = 42; Price = $null;       Volume = "  +2.86%  " },
    @{ Row = 43; Price = "2.772.91";  Volume = "  +2.67%  " },
    @{ Row = 44; Price = "0.108";     Volume = "  -0.20%  " },
    @{ Row = 45; Price = "0.255";     Volume = "  +7.06%  " },
    @{ Row = 47; Price = "123.42";    Volume = "  +0.85%  " },
    @{ Row = 48; Price = "2.06";      Volume = "  +4.69%  " },
    @{ Row = 49; Price = $null;       Volume = "  +1.56%  " },
    @{ Row = 50; Price = "34.14";     Volume = "  +21.19%  " },
    @{ Row = 51; Price = "23.81";     Volume = "  +3.13%  " }
)

# Rows whose new Price text would otherwise be auto-coerced to a number
# by plain `.Value` assignment (single decimal point, digits only) - these
# need the cell pre-set to Text format so the exact string is preserved.
$needsTextFormat = @(5, 6, 11, 14, 17, 20, 21, 22, 24, 25, 27, 30, 32, 33, 36, 37, 38, 39, 41, 44, 45, 47, 48, 50, 51)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.Price) {
        $cell = $ws.Cells.Item($r, 4)
        if ($needsTextFormat -contains $r) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.Price
    }
    $ws.Cells.Item($r, 5).Value = $u.Volume
}
